$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Packaging Overhead" (3rd sheet): selection moved to B2.
# (Done first so the workbook's active tab ends up back on sheet 4, matching
# the saved file's activeTab.)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Select()

# ---------------------------------------------------------------------------
# Sheet "Embodied Carbon" (4th sheet): insert a new column at C ("# of Parts")
# and populate DRAM / SSD rows with part details.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Mirror column B's width onto the freshly inserted column C.
$bWidth = $ws4.Columns.Item(2).ColumnWidth
$ws4.Columns.Item(3).Insert()
$ws4.Columns.Item(3).ColumnWidth = $bWidth

# New header label for the inserted column (both small tables use it).
$ws4.Cells.Item(1, 3).Value = "# of Parts"
$ws4.Cells.Item(7, 3).Value = "# of Parts"

# --- DRAM row (row 9) ---------------------------------------------------
$ws4.Cells.Item(9, 2).Value = "Samsung 8GB DDR4 PC4-21300"
$ws4.Cells.Item(9, 4).Value = 8
$ws4.Cells.Item(9, 5).Value = 10
$ws4.Cells.Item(9, 6).Value = 65
$ws4.Cells.Item(9, 7).Value = 520
$ws4.Cells.Item(9, 8).Value = "https://www.amazon.com/Samsung-PC4-21300-2666MHZ-desktop-memory/dp/B07F72RJYN"
$ws4.Cells.Item(9, 9).Value = "https://web.archive.org/web/20160706231128/http://ddr4.org/contact-us"

# --- SSD row (row 10) ----------------------------------------------------
$ws4.Cells.Item(10, 2).Value = "870 EVO SATA 2.5`" SSD 4 TB"
$ws4.Cells.Item(10, 4).Value = 4000
$ws4.Cells.Item(10, 5).Value = "N/A"
$ws4.Cells.Item(10, 6).Value = 10.7
$ws4.Cells.Item(10, 7).Value = 42800
$ws4.Cells.Item(10, 8).Value = "https://www.westerndigital.com/products/internal-drives/wd-blue-sata-2-5-ssd?ef_id=Cj0KCQiA4OybBhCzARIsAIcfn9masMsJLktL8OgSPeaWW9fMF_CxrXgOa-NTl5dSaxyt1OiJoF7ZqRcaAirvEALw_wcB:G:s&s_kwcid=AL!15012!3!!!!x!!!17824513874!&utm_medium=pdsh2&utm_source=gads&utm_campaign=WD-NA-US-PLA&utm_content=&utm_term=WDS400T2B0A#WDS400T2B0A"

# Total formula now spans the shifted columns.
$ws4.Cells.Item(13, 2).Formula = "=SUM(G9:G11)+SUM(N3:N5)"

# Selection ends up on C14 after the edits; sheet 4 remains the active tab.
$ws4.Activate()
$ws4.Range("C14").Select()
